$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New BOM rows (5-8) ---------------------------------------------------
# Row 5
$ws.Range("A5").Value = " LTC7001EMSE#PBF-ND"
$ws.Range("B5").Value = "High side driver"
$ws.Range("C5").Value = 5.13
$ws.Range("D5").Value = 1

# Row 6
$ws.Range("A6").Value = "102-4244-ND"
$ws.Range("B6").Value = "5v Buck converter"
$ws.Range("C6").Value = 2.78
$ws.Range("D6").Value = 1

# Row 7
$ws.Range("A7").Value = "1568-1233-ND"
$ws.Range("B7").Value = "Teensy"
$ws.Range("C7").Value = 14.38
$ws.Range("D7").Value = 1

# Row 8
$ws.Range("A8").Value = "565-3816-ND"
$ws.Range("B8").Value = "Aluminum capacitor"
$ws.Range("C8").Value = 4.68
$ws.Range("D8").Value = 1

# --- Totals column header --------------------------------------------------
$ws.Range("E1").Value = "Total"

# Row 9 is intentionally left blank (spacer)

# Row 10
$ws.Range("B10").Value = "Button pad"

# Row 11 (note: B entered before A to match original authoring order)
$ws.Range("B11").Value = "Supercaps"
$ws.Range("A11").Value = "1182-1018-ND"

# Rows 12-13 intentionally left blank (spacer)

# Row 14
$ws.Range("A14").Value = "TPHR6503PLL1QCT-ND"
$ws.Range("B14").Value = "Pass/short transistor"
$ws.Range("C14").Value = 2.09
$ws.Range("D14").Value = 2

# Row 15
$ws.Range("A15").Value = "1655-1928-1-ND"
$ws.Range("B15").Value = "Diode"
$ws.Range("C15").Value = 0.18
$ws.Range("D15").Value = 2

# Row 16 (note: B entered before A to match original authoring order)
$ws.Range("B16").Value = "200mOhm resistor"
$ws.Range("A16").Value = "1276-6171-1-ND"
$ws.Range("C16").Value = 0.2
$ws.Range("D16").Value = 1

# Row 17 (note: B entered before A to match original authoring order)
$ws.Range("B17").Value = "4.7k resistor"
$ws.Range("A17").Value = "RMCF0603JT4K70CT-ND"
$ws.Range("C17").Value = 0.01
$ws.Range("D17").Value = 10

# Row 18
$ws.Range("A18").Value = "587-2484-1-ND"
$ws.Range("B18").Value = "10uF cap"
$ws.Range("C18").Value = 0.2
$ws.Range("D18").Value = 5

# --- Per-line totals (column E) and grand-total in H1 ----------------------
$ws.Range("E2").Formula = "=D2*C2"
$ws.Range("E3:E28").Formula = "=D3*C3"
$ws.Range("H1").Formula = "=SUM(E2:E100)"

# --- Selection state ---------------------------------------------------
$ws.Range("C19").Select() | Out-Null
